# NATMI ligand-receptor (Wnt9a-Fzd10) sheet re-exported against a new
# TPM table. Re-write the derived numeric columns (E..T) for rows 2-4
# with the recomputed values; the identifying columns (A-D, O, P) are
# unaffected by the TPM refresh and are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3725013333333333
$ws.Range("H2").Value = 1.117504
$ws.Range("I2").Value = 0.05990362118104055
$ws.Range("J2").Value = 0.06743099450495174
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.016657
$ws.Range("N2").Value = 0.033314
$ws.Range("Q2").Value = 0.006204754709333333
$ws.Range("R2").Value = 0.037228528256
$ws.Range("S2").Value = 0.05990362118104055
$ws.Range("T2").Value = 0.06743099450495174
$ws.Range("G3").Value = 3.763367333333334
$ws.Range("I3").Value = 0.6052040917109096
$ws.Range("J3").Value = 0.6812528688240443
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.016657
$ws.Range("N3").Value = 0.033314
$ws.Range("Q3").Value = 0.06268640967133335
$ws.Range("R3").Value = 0.3761184580280001
$ws.Range("S3").Value = 0.6052040917109096
$ws.Range("T3").Value = 0.6812528688240443
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.0824755
$ws.Range("H4").Value = 4.164951
$ws.Range("I4").Value = 0.3348922871080498
$ws.Range("J4").Value = 0.2513161366710037
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.016657
$ws.Range("N4").Value = 0.033314
$ws.Range("Q4").Value = 0.03468779440350001
$ws.Range("R4").Value = 0.138751177614
$ws.Range("S4").Value = 0.3348922871080498
$ws.Range("T4").Value = 0.2513161366710037
